$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1404, 1).Value = 1403
$ws.Cells.Item(1404, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1404, 3).Value = '5:40 PM'
$ws.Cells.Item(1404, 4).Value = 'XQ421'
$ws.Cells.Item(1404, 5).Value = 'Antalya'
$ws.Cells.Item(1404, 6).Value = '(AYT)'
$ws.Cells.Item(1404, 7).Value = 'SunExpress '
$ws.Cells.Item(1404, 8).Value = 'B38M'
$ws.Cells.Item(1404, 9).Value = '(TC-SOJ)'
$ws.Cells.Item(1404, 10).Value = '6:02 PM'
$ws.Cells.Item(1404, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1404, 12).Value = '0 hours, 22 minutes'
$ws.Cells.Item(1404, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1405, 1).Value = 1404
$ws.Cells.Item(1405, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1405, 3).Value = '5:50 PM'
$ws.Cells.Item(1405, 4).Value = 'LO525'
$ws.Cells.Item(1405, 5).Value = 'Prague'
$ws.Cells.Item(1405, 6).Value = '(PRG)'
$ws.Cells.Item(1405, 7).Value = 'LOT '
$ws.Cells.Item(1405, 8).Value = 'E75S'
$ws.Cells.Item(1405, 9).Value = '(SP-LIN)'
$ws.Cells.Item(1405, 10).Value = '5:54 PM'
$ws.Cells.Item(1405, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1405, 12).Value = '0 hours, 4 minutes'
$ws.Cells.Item(1405, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1406, 1).Value = 1405
$ws.Cells.Item(1406, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1406, 3).Value = '5:55 PM'
$ws.Cells.Item(1406, 4).Value = 'LO6271'
$ws.Cells.Item(1406, 5).Value = 'Prague'
$ws.Cells.Item(1406, 6).Value = '(PRG)'
$ws.Cells.Item(1406, 7).Value = 'LOT '
$ws.Cells.Item(1406, 8).Value = 'B789'
$ws.Cells.Item(1406, 9).Value = '(SP-LSF)'
$ws.Cells.Item(1406, 10).Value = '6:11 PM'
$ws.Cells.Item(1406, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1406, 12).Value = '0 hours, 16 minutes'
$ws.Cells.Item(1406, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1407, 1).Value = 1406
$ws.Cells.Item(1407, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1407, 3).Value = '6:00 PM'
$ws.Cells.Item(1407, 4).Value = 'W61431'
$ws.Cells.Item(1407, 5).Value = 'Milan'
$ws.Cells.Item(1407, 6).Value = '(BGY)'
$ws.Cells.Item(1407, 7).Value = 'Wizz Air '
$ws.Cells.Item(1407, 8).Value = 'A321'
$ws.Cells.Item(1407, 9).Value = '(HA-LXE)'
$ws.Cells.Item(1407, 10).Value = '6:28 PM'
$ws.Cells.Item(1407, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1407, 12).Value = '0 hours, 28 minutes'
$ws.Cells.Item(1407, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1408, 1).Value = 1407
$ws.Cells.Item(1408, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1408, 3).Value = '6:20 PM'
$ws.Cells.Item(1408, 4).Value = 'W61363'
$ws.Cells.Item(1408, 5).Value = 'Basel'
$ws.Cells.Item(1408, 6).Value = '(BSL)'
$ws.Cells.Item(1408, 7).Value = 'Wizz Air '
$ws.Cells.Item(1408, 8).Value = 'A321'
$ws.Cells.Item(1408, 9).Value = '(HA-LXK)'
$ws.Cells.Item(1408, 10).Value = '6:14 PM'
$ws.Cells.Item(1408, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1408, 12).Value = '0 hours, -6 minutes'
$ws.Cells.Item(1408, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1409, 1).Value = 1408
$ws.Cells.Item(1409, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1409, 3).Value = '7:00 PM'
$ws.Cells.Item(1409, 4).Value = 'TP1207'
$ws.Cells.Item(1409, 5).Value = 'Lisbon'
$ws.Cells.Item(1409, 6).Value = '(LIS)'
$ws.Cells.Item(1409, 7).Value = 'TAP Air Portugal '
$ws.Cells.Item(1409, 8).Value = 'A320'
$ws.Cells.Item(1409, 9).Value = '(CS-TNT)'
$ws.Cells.Item(1409, 10).Value = '7:04 PM'
$ws.Cells.Item(1409, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1409, 12).Value = '0 hours, 4 minutes'
$ws.Cells.Item(1409, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1410, 1).Value = 1409
$ws.Cells.Item(1410, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1410, 3).Value = '7:25 PM'
$ws.Cells.Item(1410, 4).Value = 'AY1146'
$ws.Cells.Item(1410, 5).Value = 'Helsinki'
$ws.Cells.Item(1410, 6).Value = '(HEL)'
$ws.Cells.Item(1410, 7).Value = 'Finnair '
$ws.Cells.Item(1410, 8).Value = 'AT75'
$ws.Cells.Item(1410, 9).Value = '(OH-ATE)'
$ws.Cells.Item(1410, 10).Value = '7:36 PM'
$ws.Cells.Item(1410, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1410, 12).Value = '0 hours, 11 minutes'
$ws.Cells.Item(1410, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1411, 1).Value = 1410
$ws.Cells.Item(1411, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1411, 3).Value = '7:40 PM'
$ws.Cells.Item(1411, 4).Value = 'LPR42'
$ws.Cells.Item(1411, 5).Value = 'Katowice'
$ws.Cells.Item(1411, 6).Value = '(KTW)'
$ws.Cells.Item(1411, 7).Value = 'Polish Medical Air Rescue '
$ws.Cells.Item(1411, 8).Value = 'LJ75'
$ws.Cells.Item(1411, 9).Value = '(SP-MXS)'
$ws.Cells.Item(1411, 10).Value = '7:55 PM'
$ws.Cells.Item(1411, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1411, 12).Value = '0 hours, 15 minutes'
$ws.Cells.Item(1411, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1412, 1).Value = 1411
$ws.Cells.Item(1412, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1412, 3).Value = '7:55 PM'
$ws.Cells.Item(1412, 4).Value = 'LO137'
$ws.Cells.Item(1412, 5).Value = 'Istanbul'
$ws.Cells.Item(1412, 6).Value = '(IST)'
$ws.Cells.Item(1412, 7).Value = 'LOT '
$ws.Cells.Item(1412, 8).Value = 'E195'
$ws.Cells.Item(1412, 9).Value = '(SP-LNG)'
$ws.Cells.Item(1412, 10).Value = '8:21 PM'
$ws.Cells.Item(1412, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1412, 12).Value = '0 hours, 26 minutes'
$ws.Cells.Item(1412, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1413, 1).Value = 1412
$ws.Cells.Item(1413, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1413, 3).Value = '7:55 PM'
$ws.Cells.Item(1413, 4).Value = 'LO3825'
$ws.Cells.Item(1413, 5).Value = 'Gdansk'
$ws.Cells.Item(1413, 6).Value = '(GDN)'
$ws.Cells.Item(1413, 7).Value = 'LOT '
$ws.Cells.Item(1413, 8).Value = 'E190'
$ws.Cells.Item(1413, 9).Value = '(SP-LMB)'
$ws.Cells.Item(1413, 10).Value = '8:08 PM'
$ws.Cells.Item(1413, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1413, 12).Value = '0 hours, 13 minutes'
$ws.Cells.Item(1413, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1414, 1).Value = 1413
$ws.Cells.Item(1414, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1414, 3).Value = '7:55 PM'
$ws.Cells.Item(1414, 4).Value = 'LO3859'
$ws.Cells.Item(1414, 5).Value = 'Wroclaw'
$ws.Cells.Item(1414, 6).Value = '(WRO)'
$ws.Cells.Item(1414, 7).Value = 'LOT '
$ws.Cells.Item(1414, 8).Value = 'E75S'
$ws.Cells.Item(1414, 9).Value = '(SP-LIL)'
$ws.Cells.Item(1414, 10).Value = '8:19 PM'
$ws.Cells.Item(1414, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1414, 12).Value = '0 hours, 24 minutes'
$ws.Cells.Item(1414, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1415, 1).Value = 1414
$ws.Cells.Item(1415, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1415, 3).Value = '8:15 PM'
$ws.Cells.Item(1415, 4).Value = 'LO269'
$ws.Cells.Item(1415, 5).Value = 'Amsterdam'
$ws.Cells.Item(1415, 6).Value = '(AMS)'
$ws.Cells.Item(1415, 7).Value = 'LOT '
$ws.Cells.Item(1415, 8).Value = 'E75S'
$ws.Cells.Item(1415, 9).Value = '(SP-LIC)'
$ws.Cells.Item(1415, 10).Value = '8:28 PM'
$ws.Cells.Item(1415, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1415, 12).Value = '0 hours, 13 minutes'
$ws.Cells.Item(1415, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1416, 1).Value = 1415
$ws.Cells.Item(1416, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1416, 3).Value = '8:30 PM'
$ws.Cells.Item(1416, 4).Value = 'FZ1830'
$ws.Cells.Item(1416, 5).Value = 'Dubai'
$ws.Cells.Item(1416, 6).Value = '(DXB)'
$ws.Cells.Item(1416, 7).Value = 'flydubai '
$ws.Cells.Item(1416, 8).Value = 'B38M'
$ws.Cells.Item(1416, 9).Value = '(A6-FMM)'
$ws.Cells.Item(1416, 10).Value = '9:25 PM'
$ws.Cells.Item(1416, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1416, 12).Value = '0 hours, 55 minutes'
$ws.Cells.Item(1416, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1417, 1).Value = 1416
$ws.Cells.Item(1417, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1417, 3).Value = '9:00 PM'
$ws.Cells.Item(1417, 4).Value = 'LO789'
$ws.Cells.Item(1417, 5).Value = 'Tallinn'
$ws.Cells.Item(1417, 6).Value = '(TLL)'
$ws.Cells.Item(1417, 7).Value = 'LOT (Grzeski Livery) '
$ws.Cells.Item(1417, 8).Value = 'E195'
$ws.Cells.Item(1417, 9).Value = '(SP-LNB)'
$ws.Cells.Item(1417, 10).Value = '9:07 PM'
$ws.Cells.Item(1417, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1417, 12).Value = '0 hours, 7 minutes'
$ws.Cells.Item(1417, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1418, 1).Value = 1417
$ws.Cells.Item(1418, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1418, 3).Value = '9:00 PM'
$ws.Cells.Item(1418, 4).Value = 'W61539'
$ws.Cells.Item(1418, 5).Value = 'Reykjavik'
$ws.Cells.Item(1418, 6).Value = '(KEF)'
$ws.Cells.Item(1418, 7).Value = 'Wizz Air '
$ws.Cells.Item(1418, 8).Value = 'A21N'
$ws.Cells.Item(1418, 9).Value = '(HA-LZG)'
$ws.Cells.Item(1418, 10).Value = '9:06 PM'
$ws.Cells.Item(1418, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1418, 12).Value = '0 hours, 6 minutes'
$ws.Cells.Item(1418, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1419, 1).Value = 1418
$ws.Cells.Item(1419, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1419, 3).Value = '9:50 PM'
$ws.Cells.Item(1419, 4).Value = 'LO723'
$ws.Cells.Item(1419, 5).Value = 'Tbilisi'
$ws.Cells.Item(1419, 6).Value = '(TBS)'
$ws.Cells.Item(1419, 7).Value = 'LOT '
$ws.Cells.Item(1419, 8).Value = 'B38M'
$ws.Cells.Item(1419, 9).Value = '(SP-LVA)'
$ws.Cells.Item(1419, 10).Value = '10:04 PM'
$ws.Cells.Item(1419, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1419, 12).Value = '0 hours, 14 minutes'
$ws.Cells.Item(1419, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1420, 1).Value = 1419
$ws.Cells.Item(1420, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1420, 3).Value = '10:40 PM'
$ws.Cells.Item(1420, 4).Value = 'LO3911'
$ws.Cells.Item(1420, 5).Value = 'Krakow'
$ws.Cells.Item(1420, 6).Value = '(KRK)'
$ws.Cells.Item(1420, 7).Value = 'LOT '
$ws.Cells.Item(1420, 8).Value = 'E190'
$ws.Cells.Item(1420, 9).Value = '(SP-LMH)'
$ws.Cells.Item(1420, 10).Value = '10:48 PM'
$ws.Cells.Item(1420, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1420, 12).Value = '0 hours, 8 minutes'
$ws.Cells.Item(1420, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1421, 1).Value = 1420
$ws.Cells.Item(1421, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1421, 3).Value = '10:45 PM'
$ws.Cells.Item(1421, 4).Value = 'LO773'
$ws.Cells.Item(1421, 5).Value = 'Vilnius'
$ws.Cells.Item(1421, 6).Value = '(VNO)'
$ws.Cells.Item(1421, 7).Value = 'LOT '
$ws.Cells.Item(1421, 8).Value = 'E195'
$ws.Cells.Item(1421, 9).Value = '(SP-LNI)'
$ws.Cells.Item(1421, 10).Value = '11:15 PM'
$ws.Cells.Item(1421, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1421, 12).Value = '0 hours, 30 minutes'
$ws.Cells.Item(1421, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1422, 1).Value = 1421
$ws.Cells.Item(1422, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1422, 3).Value = '10:55 PM'
$ws.Cells.Item(1422, 4).Value = 'LO791'
$ws.Cells.Item(1422, 5).Value = 'Tallinn'
$ws.Cells.Item(1422, 6).Value = '(TLL)'
$ws.Cells.Item(1422, 7).Value = 'LOT (Retro Livery) '
$ws.Cells.Item(1422, 8).Value = 'E75S'
$ws.Cells.Item(1422, 9).Value = '(SP-LIM)'
$ws.Cells.Item(1422, 10).Value = '10:51 PM'
$ws.Cells.Item(1422, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1422, 12).Value = '0 hours, -4 minutes'
$ws.Cells.Item(1422, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1423, 1).Value = 1422
$ws.Cells.Item(1423, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1423, 3).Value = '11:00 PM'
$ws.Cells.Item(1423, 4).Value = 'LO151'
$ws.Cells.Item(1423, 5).Value = 'Tel Aviv'
$ws.Cells.Item(1423, 6).Value = '(TLV)'
$ws.Cells.Item(1423, 7).Value = 'LOT '
$ws.Cells.Item(1423, 8).Value = 'B38M'
$ws.Cells.Item(1423, 9).Value = '(SP-LVC)'
$ws.Cells.Item(1423, 10).Value = '11:20 PM'
$ws.Cells.Item(1423, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1423, 12).Value = '0 hours, 20 minutes'
$ws.Cells.Item(1423, 13).Borders.LineStyle = -4142

$ws.Cells.Item(1424, 1).Value = 1423
$ws.Cells.Item(1424, 2).Value = 'Saturday, Jan 14'
$ws.Cells.Item(1424, 3).Value = '11:05 PM'
$ws.Cells.Item(1424, 4).Value = 'LO3827'
$ws.Cells.Item(1424, 5).Value = 'Gdansk'
$ws.Cells.Item(1424, 6).Value = '(GDN)'
$ws.Cells.Item(1424, 7).Value = 'LOT '
$ws.Cells.Item(1424, 8).Value = 'E190'
$ws.Cells.Item(1424, 9).Value = '(SP-LMG)'
$ws.Cells.Item(1424, 10).Value = '11:01 PM'
$ws.Cells.Item(1424, 11).Borders.LineStyle = -4142
$ws.Cells.Item(1424, 12).Value = '0 hours, -4 minutes'
$ws.Cells.Item(1424, 13).Borders.LineStyle = -4142
